# Generate Report for Archive
#
# The en-US "Overview" sheet and its per-locale detail sheets ("zh-cn",
# "de-de") get refreshed: 01598962... and dbb5e810... move from
# "Ready for handoff" to "In Translation", and the table rows for
# dbb5e810... / bf5829eb... swap places (dbb5e810 now sits in row 4,
# bf5829eb in row 5).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A4").Value = "dbb5e810-f4b5-413a-bb77-106bd447dde6.md"
$ws.Range("B4").Value = "e2e\dbb5e810-f4b5-413a-bb77-106bd447dde6.md"
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"
$ws.Range("G4").Value = "2016-08-28 18:39:55"

$ws.Range("A5").Value = "bf5829eb-67c6-464e-84dd-39b56951552f.md"
$ws.Range("B5").Value = "e2e\bf5829eb-67c6-464e-84dd-39b56951552f.md"
$ws.Range("E5").Value = "Ready for handoff"
$ws.Range("F5").Value = "Ready for handoff"
$ws.Range("G5").Value = "2016-08-28 18:38:52"

$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"

$ws.Hyperlinks.Item(3).Name = "e2e\dbb5e810-f4b5-413a-bb77-106bd447dde6.md"
$ws.Hyperlinks.Item(4).Name = "e2e\bf5829eb-67c6-464e-84dd-39b56951552f.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "dbb5e810-f4b5-413a-bb77-106bd447dde6.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "dbb5e810-f4b5-413a-bb77-106bd447dde6.871cd8c40861f3b8ed57ca805a9fbce29b052154.zh-cn.xlf"
$ws.Range("H4").Value = "2016-08-28 18:39:51"

$ws.Range("A5").Value = "bf5829eb-67c6-464e-84dd-39b56951552f.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("G5").Value = "bf5829eb-67c6-464e-84dd-39b56951552f.89a074fef0766fe998f8cde045b10c3b6863e6c3.zh-cn.xlf"
$ws.Range("H5").Value = "2016-08-28 18:38:47"

$ws.Hyperlinks.Item(4).Name = "dbb5e810-f4b5-413a-bb77-106bd447dde6.md"
$ws.Hyperlinks.Item(5).Name = "bf5829eb-67c6-464e-84dd-39b56951552f.md"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "dbb5e810-f4b5-413a-bb77-106bd447dde6.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "dbb5e810-f4b5-413a-bb77-106bd447dde6.871cd8c40861f3b8ed57ca805a9fbce29b052154.de-de.xlf"
$ws.Range("H4").Value = "2016-08-28 18:39:55"

$ws.Range("A5").Value = "bf5829eb-67c6-464e-84dd-39b56951552f.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("G5").Value = "bf5829eb-67c6-464e-84dd-39b56951552f.89a074fef0766fe998f8cde045b10c3b6863e6c3.de-de.xlf"
$ws.Range("H5").Value = "2016-08-28 18:38:52"

$ws.Hyperlinks.Item(4).Name = "dbb5e810-f4b5-413a-bb77-106bd447dde6.md"
$ws.Hyperlinks.Item(5).Name = "bf5829eb-67c6-464e-84dd-39b56951552f.md"
